$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (column order changes: return_status now sits in col C, note in col D) ---
$ws.Range("A1").Value = "event_id"
$ws.Range("B1").Value = "customer_id"
$ws.Range("C1").Value = "return_status"
$ws.Range("D1").Value = "note"

# --- Row 2 ---
$ws.Range("A2").Value = "sn202508"
$ws.Range("B2").Value = 321891
$ws.Range("C2").Value = "TIEP_NHAN_TRA_HANG"
$ws.Range("D2").Value = "Hàng bị lỗi. Uniplo xác nhận hàng lỗi và Hương chấp nhận trả hàng"

# --- Row 3 ---
$ws.Range("A3").Value = "sn202508"
$ws.Range("B3").Value = 321903
$ws.Range("C3").Value = "DA_TRA_HANG"
$ws.Range("D3").Value = "Hương xác nhận đồng ý nhận lại hàng. 8/8/2025 - Hương đã nhận được hàng trả"

# --- Baseline formatting for the whole table: thin box border + top vertical alignment ---
$range = $ws.Range("A1:D3")
$range.Borders.LineStyle = 1
$range.Borders.Weight = 2
$range.VerticalAlignment = -4160

# --- Header row: bold + gray fill ---
$header = $ws.Range("A1:D1")
$header.Font.Bold = $true
$header.Interior.Color = 12566463

# --- B2 (321891): no border, explicit black font color, keep top alignment ---
$b2 = $ws.Range("B2")
$b2.Borders.LineStyle = 0
$b2.Font.Color = 0
$b2.VerticalAlignment = -4160

# --- D3: wrap text, taller row ---
$ws.Range("D3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 30.75

# --- Selection moved to B8 ---
$ws.Range("B8").Select()
